$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new header for column C
$ws.Range("C1").Value = "data_size"

# Fill in the data_size values for each topic row
$ws.Range("C2").Value = 48
$ws.Range("C3").Value = 282
$ws.Range("C4").Value = 132
$ws.Range("C5").Value = 288
$ws.Range("C6").Value = 41
$ws.Range("C7").Value = 77
$ws.Range("C8").Value = 17
$ws.Range("C9").Value = 23
$ws.Range("C10").Value = 22
$ws.Range("C11").Value = 4
$ws.Range("C12").Value = 82

# Update the active selection to match the target state
$ws.Range("Q10").Select()
